$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# Status columns for zh-cn (E) and de-de (F), plus the "Latest HO Xliff
# Generate Date" column (G) move from "In Translation" to "Ready for
# handoff", and the generation timestamp ticks forward by 40s.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 14:42:40"

# --- zh-cn sheet ------------------------------------------------------
# Status (C) and Latest Handoff Datetime (H)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 14:42:36"

# --- de-de sheet --------------------------------------------------------
# Status (C) and Latest Handoff Datetime (H)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-04 14:42:40"

# --- Column widths -----------------------------------------------------
# The "Ready for handoff" status text is longer than "In Translation", so
# the status columns are widened to fit: Overview!E:F and the "Status"
# column (C) on both the zh-cn and de-de sheets.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
